$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 7 (existing rows 7:31 shift down to 8:32),
# inheriting formatting from the row above (row 6), same as Excel's
# default "Insert" behavior.
$ws.Rows("7:7").Insert()

# Populate the newly inserted row 7 with the new item "2:Job" - a
# job-scoped variant of item "2" (now row 9 after the insert), carrying
# the same description/account/price/income account/type.
$ws.Range("B7").Value = "2:Job"
$ws.Range("C7").Value = $ws.Range("C9").Value2
$ws.Range("D7").Value = $ws.Range("D9").Value2
$ws.Range("E7").Value = $ws.Range("E9").Value2
$ws.Range("G7").Value = $ws.Range("G9").Value2
$ws.Range("H7").Value = $ws.Range("H9").Value2
$ws.Range("J7").Value = $ws.Range("J9").Value2

# Update the saved selection on the sheet.
$ws.Range("P7").Select()
